$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "42.925.37"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "2.571.74"
$ws.Range("E3").Value = "  +3.84%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "302.87"
$ws.Range("E5").Value = "  +3.86%  "
$ws.Range("D6").Value = "97.19"
$ws.Range("E6").Value = "  +6.43%  "
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").Value = "36.87"
$ws.Range("E10").Value = "  +4.07%  "
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  +7.60%  "
$ws.Range("D14").Value = "2.578.55"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "0.884"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").Value = "14.35"
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("D17").Value = "42.977.75"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "13.04"
$ws.Range("E18").Value = "  +8.77%  "
$ws.Range("D19").Value = "0.0₃0995"
$ws.Range("E19").Value = "  +5.50%  "
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +4.14%  "
$ws.Range("D21").Value = "71.99"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "254.81"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "2.97"
$ws.Range("E23").Value = "  +5.18%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "28.62"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "10.25"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("D28").Value = "37.81"
$ws.Range("E28").Value = "  +5.24%  "
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "155.60"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").Value = "2.18"
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("D33").Value = "2.75"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "0.0812"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("D36").Value = "18.43"
$ws.Range("E36").Value = "  +13.15%  "
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "23.52"
$ws.Range("E39").Value = "  -13.28%  "
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "3.88"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0310"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("D43").Value = "2.06"
$ws.Range("E43").Value = "  +27.86%  "
$ws.Range("D44").Value = "2.065.39"
$ws.Range("E44").Value = "  +3.67%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "9.29"
$ws.Range("E46").Value = "  +5.84%  "
$ws.Range("D47").Value = "85.24"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").Value = "77.49"
$ws.Range("E48").Value = "  +17.19%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.823.00"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "106.19"
$ws.Range("E50").Value = "  +5.33%  "
$ws.Range("E51").Value = "  +4.80%  "
